$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 45212 to 45221
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
